$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

# Type the bold "Meta description" label first
$metaRange = $metaPara.Range
$metaRange.Collapse(1)
$metaRange.InsertAfter("Meta description")

# Bold just the label we typed
$labelRange = $d.Range($metaPara.Range.Start, $metaPara.Range.Start + 17)
$labelRange.Font.Bold = $true

# Now append the (non-bold) description text right before the paragraph mark
$metaPara2 = $d.Paragraphs.Item(2)
$insPos = $metaPara2.Range.End - 1
$descText = ": Experience an immersive adventure - Play 50 Dragons for free and enjoy stunning graphics, bonus features and a potential win of up to €125,000."
$descInsertionPoint = $d.Range($insPos, $insPos)
$descInsertionPoint.InsertAfter($descText)

$descRange = $d.Range($insPos, $insPos + $descText.Length)
$descRange.Font.Bold = $false

# ------------------------------------------------------------------
# 2) Remove the duplicate bold title paragraph near the end of the doc
#    ("Play 50 Dragons Slot Free - Mesmerizing Design with Lucrative Wins")
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($count - 1)
$dupTitlePara.Range.Delete()

# ------------------------------------------------------------------
# 3) Replace the final (italic) paragraph's text with the image prompt.
#    We overwrite just the visible text (Range excluding the trailing
#    paragraph mark) via direct Range.Text assignment rather than
#    Find/Replace, since Find's replacement path "smart-quotes" every
#    apostrophe it inserts - direct Range.Text keeps the characters we
#    typed exactly as-is (including the one straight apostrophe in the
#    source text) while the run's existing italic formatting carries
#    over unchanged because we never touch the run's rPr.
# ------------------------------------------------------------------
$newText = "Prompt: Create a feature image for the online slot game “50 Dragons.” The image should feature a cartoon-style happy Maya warrior with glasses. The warrior should be surrounded by sparkling gold dragons, tigers, masks and peacocks. Use gold as the dominant color scheme to add to the game's luxurious and exotic feel. Make the Maya warrior look adventurous and ready for a thrilling game in the mystical world of the Orient. The image should represent the game’s stunning design and exciting features, such as the Pearl symbol and Scatter symbol that trigger Free Spins bonus for even bigger rewards. The image should encourage new players to experience the adventure and big rewards of 50 Dragons by Aristocrat."

$finalCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($finalCount)
$lastRange = $lastPara.Range
$textRange = $d.Range($lastRange.Start, $lastRange.End - 1)
$textRange.Text = $newText

Write-Host "Final paragraph count:" $d.Paragraphs.Count
